$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# The sheet used to store a stray duplicate-of-row-2 "header" row and only
# 7 columns (B:G). Bring it in line with the other property sheets:
#  - row 1 becomes real column headers
#  - rows 2/3 gain the common property_category..index metadata columns (H:N)

# Mirror the existing header/data formatting onto the new columns first so
# the new cells inherit the same bold/border (row1) / plain (rows2-3) style
# that the rest of the row already uses.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 1: header labels ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Rows 2 & 3: fill in the new metadata columns (H:N) ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("K2").Value = "田秋堇"
$ws.Range("L2").Value = 1316
$ws.Range("M2").Value = "tmp9b251"
$ws.Range("N2").Value = 30

$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("K3").Value = "田秋堇"
$ws.Range("L3").Value = 1316
$ws.Range("M3").Value = "tmp9b251"
$ws.Range("N3").Value = 31

# The "date" column (J) holds literal text "2012-04-10" -- assigning that
# string directly gets auto-converted to a real Excel date serial, so route
# it through a formula + paste-as-values so it lands back as literal text.
$ws.Range("J2").Formula = '="2012-04-10"'
$ws.Range("J2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

$ws.Range("J3").Formula = '="2012-04-10"'
$ws.Range("J3").Copy()
$ws.Range("J3").PasteSpecial(-4163)

$excel.CutCopyMode = 0
